$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
# Set A2's date-like text without letting Excel auto-convert it into a date
# serial: build it as a formula that evaluates to the literal text, then
# flatten the formula down to a static value via copy/paste-special, which
# preserves the cell's existing style (unlike touching NumberFormat).
$ws.Range("A2").Formula = "=""15-JAN-26"""
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4163) | Out-Null  # xlPasteValues

$ws.Range("D2").Value = 10295
$ws.Range("E2").Value = 10660
$ws.Range("F2").Value = -365

# --- Row 3 updates ---
$ws.Range("D3").Value = 5998
$ws.Range("E3").Value = 6212
$ws.Range("F3").Value = -214

# --- Remove row 4 entirely (shrinks used range / dimension to A1:K3) ---
$ws.Rows(4).Delete()

$excel.CutCopyMode = $false
